$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 22:31"

# Row data: Row number, Country name, then values for columns B..H
$rows = @(
    @{ Row = 4; Name = "Estados Unidos"; Vals = @(2710176, 28365, 1126719, 1454159, 0, 515, 129298) },
    @{ Row = 17; Name = "Alemania"; Vals = @(195781, 389, 179100, 7630, 0, 10, 9051) },
    @{ Row = 19; Name = "Francia"; Vals = @(164801, 541, 76274, 58684, 0, 30, 29843) },
    @{ Row = 20; Name = "Sudafrica"; Vals = @(151209, 6945, 73543, 75009, 0, 128, 2657) },
    @{ Row = 21; Name = "Banglades"; Vals = @(145483, 3682, 59624, 84012, 0, 64, 1847) },
    @{ Row = 26; Name = "Egipto"; Vals = @(68311, 1557, 18460, 46898, 0, 81, 2953) },
    @{ Row = 27; Name = "Suecia"; Vals = @(67667, 0, 0, 0, 0, 0, 5310) },
    @{ Row = 52; Name = "Israel"; Vals = @(25244, 803, 17341, 7583, 0, 1, 320) },
    @{ Row = 53; Name = "Nigeria"; Vals = @(25133, 0, 9402, 15158, 0, 0, 573) },
    @{ Row = 70; Name = "Costa de Marfil"; Vals = @(9499, 285, 4273, 5158, 0, 2, 68) },
    @{ Row = 71; Name = "Sudan"; Vals = @(9257, 0, 4014, 4671, 0, 0, 572) },
    @{ Row = 87; Name = "Guinea"; Vals = @(5391, 40, 4326, 1032, 0, 2, 33) },
    @{ Row = 96; Name = "Republica de Africa Central"; Vals = @(3745, 132, 787, 2911, 0, 0, 47) },
    @{ Row = 97; Name = "Costa Rica"; Vals = @(3459, 190, 1436, 2007, 0, 1, 16) },
    @{ Row = 98; Name = "Grecia"; Vals = @(3409, 19, 1374, 1843, 0, 1, 192) },
    @{ Row = 100; Name = "Somalia"; Vals = @(2924, 20, 910, 1924, 0, 0, 90) },
    @{ Row = 104; Name = "Estado de Palestina"; Vals = @(2428, 243, 451, 1970, 0, 2, 7) },
    @{ Row = 121; Name = "Zambia"; Vals = @(1594, 26, 1329, 241, 0, 2, 24) },
    @{ Row = 152; Name = "Zimbabue"; Vals = @(591, 17, 162, 422, 0, 0, 7) },
    @{ Row = 161; Name = "Comoras"; Vals = @(303, 31, 200, 96, 0, 0, 7) },
    @{ Row = 162; Name = "Birmania"; Vals = @(299, 0, 222, 71, 0, 0, 6) },
    @{ Row = 163; Name = "Angola"; Vals = @(284, 8, 93, 178, 0, 2, 13) },
    @{ Row = 164; Name = "Siria"; Vals = @(279, 10, 105, 165, 0, 0, 9) },
    @{ Row = 169; Name = "Namibia"; Vals = @(205, 9, 24, 181, 0, 0, 0) },
    @{ Row = 170; Name = "Eritrea"; Vals = @(203, 12, 53, 150, 0, 0, 0) },
    @{ Row = 171; Name = "Islas Caimanes"; Vals = @(199, 0, 189, 9, 0, 0, 1) },
    @{ Row = 203; Name = "Laos"; Vals = @(19, 0, 19, 0, 0, 0, 0) },
    @{ Row = 204; Name = "Santa Lucia"; Vals = @(19, 0, 19, 0, 0, 0, 0) },
    @{ Row = 205; Name = "Fiyi"; Vals = @(18, 0, 18, 0, 0, 0, 0) },
    @{ Row = 206; Name = "Dominica"; Vals = @(18, 0, 18, 0, 0, 0, 0) }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    for ($i = 0; $i -lt $r.Vals.Length; $i++) {
        $ws.Cells.Item($r.Row, $i + 2).Value = $r.Vals[$i]
    }
}

Write-Host "Updated $($rows.Length) rows"
